$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update poll figures for this round (Essential / ResolvePM).
$ws.Range("A2").Value = 36
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 1

# Add the next poll's raw numbers below, ready for the next normalisation pass.
$ws.Range("A8").Value = "LNP"
$ws.Range("B8").Value = "ALP"
$ws.Range("D8").Formula = "=SUM(A9:B9)"
$ws.Range("A9").Value = 48
$ws.Range("B9").Value = 45

# Leave selection where Excel would land after entering the last value.
$ws.Range("A10").Select()
